# Updates the cryptos list (Sheet1) with refreshed price/volume data.
# For numeric-looking Price values (column D) we briefly force a text
# number format before assigning the literal string, then restore the
# "Normal" style so the cell keeps its original (default) style index -
# this stops Excel from silently re-parsing strings like "246.85" as
# actual numbers while avoiding any spurious style diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.532.06'
$ws.Range("E2").Value = '  +1.30%  '

$ws.Range("D3").Value = '1.909.53'
$ws.Range("E3").Value = '  +2.82%  '

$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.660'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.05%  '

$ws.Range("E7").Value = '  +0.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.05'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.36%  '

$ws.Range("E9").Value = '  +4.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0716'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.93%  '

$ws.Range("E12").Value = '  +1.12%  '

$ws.Range("D13").Value = '2.184.74'
$ws.Range("E13").Value = '  +2.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.11%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.701'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.38%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.915.85'
$ws.Range("E16").Value = '  +2.60%  '

$ws.Range("E17").Value = '  +3.52%  '

$ws.Range("D18").Value = '35.528.14'
$ws.Range("E18").Value = '  +1.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.31%  '

$ws.Range("D20").Value = '0.0₃0832'
$ws.Range("E20").Value = '  +4.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '244.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.67'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.38%  '

$ws.Range("E23").Value = '  +2.51%  '

$ws.Range("E25").Value = '  +1.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +15.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '171.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.86%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.129'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.964'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +22.79%  '

$ws.Range("E33").Value = '  +2.19%  '

$ws.Range("E34").Value = '  +5.77%  '

$ws.Range("E35").Value = '  +0.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.72'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.14%  '

$ws.Range("E37").Value = '  +0.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.38%  '

$ws.Range("E39").Value = '  +2.67%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0205'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.18%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '92.28'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.51%  '

$ws.Range("E42").Value = '  +17.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.25%  '

$ws.Range("D44").Value = '1.348.91'
$ws.Range("E44").Value = '  -0.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +38.95%  '

$ws.Range("B47").Value = 'Gas'
$ws.Range("C47").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.64%  '

$ws.Range("B48").Value = 'HuobiToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.46%  '

$ws.Range("D51").Value = '2.095.08'
$ws.Range("E51").Value = '  +2.67%  '
